# disk_savvy.xlsx - append newer disk-usage sampling rows (40-57) captured
# by the scraping script, mirroring the formatting of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date (col A), time-of-day (col B), file count (col C), disk_space MB (col D)
$newRows = @(
    @(44697,0.63185185185185189,36208,685.2),
    @(44697,0.6348611111111111,36201,685.21),
    @(44697,0.96414351851851843,36204,685.22),
    @(44698,0.40944444444444444,36252,685.72),
    @(44698,0.99490740740740735,36250,684.67),
    @(44699,0.93873842592592593,36254,684.85),
    @(44700,0.48609953703703707,36248,684.7),
    @(44700,0.98951388888888892,36248,684.81),
    @(44701,0.99653935185185183,36250,684.78),
    @(44704,0.64586805555555549,36250,684.79),
    @(44706,0.56009259259259259,36250,684.8),
    @(44712,0.47394675925925928,36253,684.81),
    @(44713,0.65696759259259252,36259,684.98),
    @(44715,0.53877314814814814,36259,685.1),
    @(44723,0.46240740740740738,36231,685.14),
    @(44766,0.85538194444444438,36354,685.55),
    @(44769,0.93134259259259267,70024,1140),
    @(44771,0.76491898148148152,71693,1190)
)

$startRow = 40
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Seed the new row with the previous row's formatting (date / time
    # number formats on columns A and B) before writing the values.
    $ws.Range("A" + ($r - 1) + ":D" + ($r - 1)).Copy()
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$lastRow = $startRow + $newRows.Count - 1

# Move the active selection/view to just past the newly appended data,
# matching where Excel would land after entering the last row.
[void]$ws.Range("A" + $lastRow).Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A" + ($lastRow + 1)).Select()
